{"js": "// Update the division expressions in the worksheet table.\n// Each cell holds a single run of text like \"295\u00f74=\" which must become\n// \"320\u00f76=\", etc. All old values are unique in the document, so a simple\n// exact-text search+replace per pair is safe and preserves run formatting\n// (font/size) because we only rewrite the text of the matched range.\nconst replacements = [\n  [\"295\u00f74=\", \"320\u00f76=\"],\n  [\"234\u00f77=\", \"417\u00f79=\"],\n  [\"939\u00f79=\", \"723\u00f78=\"],\n  [\"204\u00f73=\", \"310\u00f73=\"],\n  [\"587\u00f74=\", \"947\u00f74=\"],\n  [\"291\u00f74=\", \"273\u00f74=\"],\n  [\"517\u00f79=\", \"885\u00f78=\"],\n  [\"317\u00f74=\", \"127\u00f74=\"],\n  [\"130\u00f75=\", \"874\u00f78=\"],\n  [\"745\u00f78=\", \"666\u00f75=\"],\n  [\"454\u00f75=\", \"419\u00f78=\"],\n  [\"406\u00f73=\", \"732\u00f78=\"],\n  [\"163\u00f78=\", \"383\u00f74=\"],\n  [\"851\u00f75=\", \"762\u00f74=\"],\n  [\"894\u00f79=\", \"479\u00f75=\"],\n  [\"445\u00f74=\", \"789\u00f78=\"],\n  [\"501\u00f73=\", \"316\u00f74=\"],\n  [\"473\u00f77=\", \"475\u00f76=\"],\n  [\"868\u00f74=\", \"929\u00f78=\"],\n  [\"623\u00f78=\", \"562\u00f72=\"],\n  [\"172\u00f74=\", \"270\u00f76=\"],\n  [\"482\u00f72=\", \"480\u00f78=\"],\n  [\"963\u00f76=\", \"950\u00f77=\"],\n  [\"896\u00f78=\", \"229\u00f72=\"],\n  [\"855\u00f74=\", \"598\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division expressions in the worksheet table.\n# Each cell holds a single run of text like \"295\u00f74=\" which must become\n# \"320\u00f76=\", etc. All old values are unique in the document, so exact-text\n# Find/Replace per pair is safe and leaves run formatting (font/size)\n# untouched since only the matched text is replaced.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"295\u00f74=\", \"320\u00f76=\"),\n    @(\"234\u00f77=\", \"417\u00f79=\"),\n    @(\"939\u00f79=\", \"723\u00f78=\"),\n    @(\"204\u00f73=\", \"310\u00f73=\"),\n    @(\"587\u00f74=\", \"947\u00f74=\"),\n    @(\"291\u00f74=\", \"273\u00f74=\"),\n    @(\"517\u00f79=\", \"885\u00f78=\"),\n    @(\"317\u00f74=\", \"127\u00f74=\"),\n    @(\"130\u00f75=\", \"874\u00f78=\"),\n    @(\"745\u00f78=\", \"666\u00f75=\"),\n    @(\"454\u00f75=\", \"419\u00f78=\"),\n    @(\"406\u00f73=\", \"732\u00f78=\"),\n    @(\"163\u00f78=\", \"383\u00f74=\"),\n    @(\"851\u00f75=\", \"762\u00f74=\"),\n    @(\"894\u00f79=\", \"479\u00f75=\"),\n    @(\"445\u00f74=\", \"789\u00f78=\"),\n    @(\"501\u00f73=\", \"316\u00f74=\"),\n    @(\"473\u00f77=\", \"475\u00f76=\"),\n    @(\"868\u00f74=\", \"929\u00f78=\"),\n    @(\"623\u00f78=\", \"562\u00f72=\"),\n    @(\"172\u00f74=\", \"270\u00f76=\"),\n    @(\"482\u00f72=\", \"480\u00f78=\"),\n    @(\"963\u00f76=\", \"950\u00f77=\"),\n    @(\"896\u00f78=\", \"229\u00f72=\"),\n    @(\"855\u00f74=\", \"598\u00f75=\")\n)\n\n# WdFindWrap.wdFindContinue = 1, WdReplace.wdReplaceAll = 2 (the interpreter\n# does not pre-seed the Wd* enum names as variables, so literal values are\n# used directly in the Execute() call below).\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n"}
